$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: Volume/Number and report week dates
$ws.Range("A8").Value = "Volume 31   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/4/2024  Through  11/10/2024"

# Weekly crime statistics data updates (rows 15-31)
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 20
$ws.Range("K15").Value = 33.333333333333
$ws.Range("L15").Value = -9.090909090909
$ws.Range("M15").Value = 11.111111111111
$ws.Range("N15").Value = -37.5
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 7.692307692307
$ws.Range("I16").Value = 165
$ws.Range("J16").Value = 137
$ws.Range("K16").Value = 20.437956204379
$ws.Range("L16").Value = -4.069767441860
$ws.Range("M16").Value = -21.800947867298
$ws.Range("N16").Value = -81.808158765159
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 26.086956521739
$ws.Range("I17").Value = 286
$ws.Range("J17").Value = 252
$ws.Range("K17").Value = 13.492063492063
$ws.Range("L17").Value = 7.116104868913
$ws.Range("M17").Value = 74.390243902439
$ws.Range("N17").Value = -37.690631808278
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 147
$ws.Range("J18").Value = 116
$ws.Range("K18").Value = 26.724137931034
$ws.Range("L18").Value = -33.484162895927
$ws.Range("M18").Value = -33.484162895927
$ws.Range("N18").Value = -85.431119920713
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -45.454545454545
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -27.083333333333
$ws.Range("I19").Value = 413
$ws.Range("J19").Value = 490
$ws.Range("K19").Value = -15.714285714285
$ws.Range("L19").Value = -42.075736325385
$ws.Range("M19").Value = 28.260869565217
$ws.Range("N19").Value = 2.736318407960
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = -45
$ws.Range("I20").Value = 133
$ws.Range("J20").Value = 159
$ws.Range("K20").Value = -16.352201257861
$ws.Range("L20").Value = 4.724409448818
$ws.Range("M20").Value = 23.148148148148
$ws.Range("N20").Value = -84.334511189634
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -28.571428571428
$ws.Range("F21").Value = 100
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = -12.280701754386
$ws.Range("I21").Value = 1167
$ws.Range("J21").Value = 1172
$ws.Range("K21").Value = -0.426621160409
$ws.Range("L21").Value = -23.475409836065
$ws.Range("M21").Value = 11.142857142857
$ws.Range("N21").Value = -68.244897959183
$ws.Range("D22").Value = 3
$ws.Range("J22").Value = 32
$ws.Range("K22").Value = -40.625
$ws.Range("M22").Value = -26.923076923076
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 83
$ws.Range("H24").Value = -3.488372093023
$ws.Range("I24").Value = 842
$ws.Range("J24").Value = 1085
$ws.Range("K24").Value = -22.396313364055
$ws.Range("L24").Value = -29.005059021922
$ws.Range("M24").Value = 23.099415204678
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -8.333333333333
$ws.Range("G25").Value = 47
$ws.Range("H25").Value = -19.148936170212
$ws.Range("I25").Value = 342
$ws.Range("J25").Value = 612
$ws.Range("K25").Value = -44.117647058823
$ws.Range("L25").Value = -48.103186646434
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = 58.333333333333
$ws.Range("F26").Value = 56
$ws.Range("G26").Value = 45
$ws.Range("H26").Value = 24.444444444444
$ws.Range("I26").Value = 492
$ws.Range("J26").Value = 461
$ws.Range("K26").Value = 6.724511930585
$ws.Range("L26").Value = 27.461139896373
$ws.Range("M26").Value = -10.869565217391
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 27
$ws.Range("K27").Value = 22.727272727272
$ws.Range("L27").Value = -12.903225806451
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 50
$ws.Range("J28").Value = 73
$ws.Range("K28").Value = -31.506849315068
$ws.Range("L28").Value = -12.280701754386
$ws.Range("J31").Value = 4
$ws.Range("K31").Value = 75
